$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (data rows 2-51) to Text format so that numeric-looking
# strings (e.g. "291.07") are stored verbatim as text, matching the source data,
# instead of being auto-converted into floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "39.784.31"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "2.206.19"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "291.07"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "86.51"
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.466"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "30.24"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "50.00"
$ws.Range("E11").Value = "  +6.01%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").Value = "0.0776"
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("E13").Value = "  +2.59%  "
$ws.Range("D14").Value = "6.40"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").Value = "2.557.91"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "13.67"
$ws.Range("E16").Value = "  -2.72%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.175.30"
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("D18").Value = "0.727"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "39.749.54"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "0.0₃0881"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "11.15"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").Value = "5.72"
$ws.Range("E22").Value = "  -1.49%  "
$ws.Range("D23").Value = "65.36"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "235.99"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "2.44"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").Value = "1.82"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").Value = "23.27"
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("E29").Value = "  -2.90%  "
$ws.Range("D30").Value = "9.19"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("D31").Value = "157.16"
$ws.Range("E31").Value = "  +3.57%  "
$ws.Range("D32").Value = "31.70"
$ws.Range("E32").Value = "  -3.57%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").Value = "4.93"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "0.0707"
$ws.Range("E35").Value = "  -1.44%  "
$ws.Range("D36").Value = "2.89"
$ws.Range("E36").Value = "  +3.15%  "
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").Value = "0.0978"
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").Value = "15.13"
$ws.Range("E41").Value = "  -4.12%  "
$ws.Range("D42").Value = "2.108.06"
$ws.Range("E42").Value = "  +2.10%  "
$ws.Range("D43").Value = "3.70"
$ws.Range("E43").Value = "  -2.37%  "
$ws.Range("D44").Value = "0.0267"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "2.10"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "9.88"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "17.68"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").Value = "2.69"
$ws.Range("E48").Value = "  +2.84%  "
$ws.Range("D49").Value = "2.431.18"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").Value = "1.47"
$ws.Range("E50").Value = "  +2.90%  "
$ws.Range("D51").Value = "88.09"
$ws.Range("E51").Value = "  -1.06%  "
